# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-250) from 2023-09-06 (serial 45175) to 2023-09-08 (serial 45177).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C250").Value = 45177
